# Applies the "Abstract.docx" diff:
#   1) Removes the trailing, text-less <w:r> run (an <w:rPr><w:rtl .../></w:rPr>
#      run with no <w:t>) that follows each of three bullet paragraphs:
#        - "Functional: Monads wrapped Values."
#        - "Functional: Monads wrapped Value Types."
#        - "Functional: Monads Transforms."
#   2) Expands the text of the "Functional Protocol: Dataflow / Parsing. ..."
#      bullet and drops its paragraph-mark underline override
#      (<w:rPr><w:u w:val="none"/></w:rPr> -> no override).

$d = $word.ActiveDocument

$wdNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
# Shared bullet pPr (numPr ilvl=0/numId=2, ind left=600 hanging=360) used by
# every item in this list - unchanged by the diff.
$pPr   = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:left="600" w:hanging="360"/></w:pPr>'

function Find-ParagraphByPrefix($prefix) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

function Replace-ParagraphSingleRun($prefix, $text) {
    $para = Find-ParagraphByPrefix($prefix)
    if ($null -eq $para) {
        throw "Paragraph starting with '$prefix' not found"
    }
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xml = "<w:p $wdNs>$pPr<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">$escaped</w:t></w:r></w:p>"
    $para.Range.InsertXML($xml)
}

# 1) Drop the stray trailing empty run on these three bullets (text is
#    otherwise unchanged).
Replace-ParagraphSingleRun "Functional: Monads wrapped Values." "Functional: Monads wrapped Values."
Replace-ParagraphSingleRun "Functional: Monads wrapped Value Types." "Functional: Monads wrapped Value Types."
Replace-ParagraphSingleRun "Functional: Monads Transforms." "Functional: Monads Transforms."

# 2) Expand the "Functional Protocol" bullet's text and remove the
#    paragraph-mark underline override at the same time.
Replace-ParagraphSingleRun "Functional Protocol: Dataflow" "Functional Protocol: Dataflow / Parsing. Monads Wrapper Types: MVC / DCI Graph (Node, Type / Context, Arc, Node); Discrete / Continuous Models CSPOs. Model driven (Types Domain / Range ordered contexts) Transforms (Augmentations)."
